$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 163.25
$ws.Range("I2").Value = 93.85714
$ws.Range("K2").Value = 93.85714
$ws.Range("M2").Value = 19.14286

$ws.Range("H43").Value = 4235.8
$ws.Range("I43").Value = 2373.8
$ws.Range("J43").Value = 6097.8
$ws.Range("K43").Value = 2373.8
$ws.Range("L43").Value = 6097.8
$ws.Range("M43").Value = -2304.8
$ws.Range("N43").Value = -6235.8

$ws.Range("H62").Value = 12828003
$ws.Range("J62").Value = 8878
$ws.Range("L62").Value = 8878
$ws.Range("N62").Value = -10126

$ws.Range("H65").Value = 12828003
$ws.Range("J65").Value = 8878
$ws.Range("L65").Value = 44390
$ws.Range("N65").Value = -50630

$ws.Range("H86").Value = 1399.5
$ws.Range("I86").Value = 1399.5
$ws.Range("K86").Value = 1399.5
$ws.Range("M86").Value = -276.5

$ws.Range("H89").Value = 1399.5
$ws.Range("I89").Value = 1399.5
$ws.Range("K89").Value = 6997.5
$ws.Range("M89").Value = -1381.5

$ws.Range("H106").Value = 7767.625
$ws.Range("I106").Value = 1734.4286
$ws.Range("K106").Value = 1734.4286
$ws.Range("M106").Value = -1103.4286

$ws.Range("H132").Value = 3795.7917
$ws.Range("I132").Value = 3676.3809
$ws.Range("K132").Value = 11029.1427
$ws.Range("M132").Value = -8499.1427

$ws.Range("H137").Value = 20028.096
$ws.Range("I137").Value = 28706.023
$ws.Range("J137").Value = 3855.5908
$ws.Range("K137").Value = 86118.069
$ws.Range("L137").Value = 11566.7724
$ws.Range("M137").Value = -83568.069
$ws.Range("N137").Value = -16666.7724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2358.2222
$ws.Range("I63").Value = 2358.2222
$ws.Range("K63").Value = 2358.2222
$ws.Range("M63").Value = -1672.2222

$ws.Range("H66").Value = 2358.2222
$ws.Range("I66").Value = 2358.2222
$ws.Range("K66").Value = 11791.111
$ws.Range("M66").Value = -8359.111000000001

$ws.Range("H74").Value = 183281.81
$ws.Range("I74").Value = 222900.56
$ws.Range("K74").Value = 222900.56
$ws.Range("M74").Value = -222026.56

$ws.Range("H77").Value = 183281.81
$ws.Range("I77").Value = 222900.56
$ws.Range("K77").Value = 1114502.8
$ws.Range("M77").Value = -1110134.8

$ws.Range("H88").Value = 2329.5833
$ws.Range("I88").Value = 2363.8667
$ws.Range("K88").Value = 2363.8667
$ws.Range("M88").Value = -1957.8667

$ws.Range("H91").Value = 2329.5833
$ws.Range("I91").Value = 2363.8667
$ws.Range("K91").Value = 2363.8667
$ws.Range("M91").Value = -959.8667

$ws.Range("H97").Value = 1669.6364
$ws.Range("I97").Value = 1349
$ws.Range("K97").Value = 1349
$ws.Range("M97").Value = -853

$ws.Range("H132").Value = 2595.257
$ws.Range("I132").Value = 2844.4167
$ws.Range("J132").Value = 2465.261
$ws.Range("K132").Value = 8533.250100000001
$ws.Range("L132").Value = 7395.782999999999
$ws.Range("M132").Value = -6003.250100000001
$ws.Range("N132").Value = -12455.783

$ws.Range("H134").Value = 150428.5
$ws.Range("J134").Value = 150428.5
$ws.Range("L134").Value = 150428.5
$ws.Range("N134").Value = -160568.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5350.207
$ws.Range("I134").Value = 5350.207
$ws.Range("K134").Value = 16050.621
$ws.Range("M134").Value = -13515.621

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1784.7142
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 1784.7142
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 1784.7142
$ws.Range("N10").Value = -2062.7142
$ws.Range("M10").ClearContents()

$ws.Range("H31").Value = 528654.5
$ws.Range("J31").Value = 2832.3333
$ws.Range("L31").Value = 2832.3333
$ws.Range("N31").Value = -3422.3333

$ws.Range("H34").Value = 528654.5
$ws.Range("J34").Value = 2832.3333
$ws.Range("L34").Value = 2832.3333
$ws.Range("N34").Value = -3236.3333

$ws.Range("H132").Value = 14489.4
$ws.Range("I132").Value = 12128.571
$ws.Range("J132").Value = 19998
$ws.Range("K132").Value = 36385.713
$ws.Range("L132").Value = 59994
$ws.Range("M132").Value = -33855.713
$ws.Range("N132").Value = -65054

$ws.Range("H134").Value = 3518.2
$ws.Range("I134").Value = 3518.2
$ws.Range("K134").Value = 10554.6
$ws.Range("M134").Value = -8019.599999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1359.8
$ws.Range("I5").Value = 949.5
$ws.Range("K5").Value = 2848.5
$ws.Range("M5").Value = -2736.5

$ws.Range("H12").Value = 4628.125
$ws.Range("I12").Value = 432.83334
$ws.Range("J12").Value = 7145.3
$ws.Range("K12").Value = 1298.50002
$ws.Range("L12").Value = 21435.9
$ws.Range("M12").Value = -1125.50002
$ws.Range("N12").Value = -21781.9

$ws.Range("H38").Value = 265.42856
$ws.Range("I38").Value = 307.5
$ws.Range("J38").Value = 237.38095
$ws.Range("K38").Value = 922.5
$ws.Range("L38").Value = 712.1428500000001
$ws.Range("M38").Value = -575.5
$ws.Range("N38").Value = -1406.14285

$ws.Range("H68").Value = 13889841
$ws.Range("J68").Value = 1324.25
$ws.Range("L68").Value = 3972.75
$ws.Range("N68").Value = -5594.75

$ws.Range("H71").Value = 13889841
$ws.Range("J71").Value = 1324.25
$ws.Range("L71").Value = 11918.25
$ws.Range("N71").Value = -20030.25

$ws.Range("H113").Value = 713.1667
$ws.Range("I113").Value = 535
$ws.Range("K113").Value = 1605
$ws.Range("M113").Value = 565

$ws.Range("H121").Value = 2200.158
$ws.Range("I121").Value = 2152.6365
$ws.Range("J121").Value = 2265.5
$ws.Range("K121").Value = 6457.9095
$ws.Range("L121").Value = 6796.5
$ws.Range("M121").Value = -5147.9095
$ws.Range("N121").Value = -9416.5

$ws.Range("H135").Value = 1359.8
$ws.Range("I135").Value = 949.5
$ws.Range("K135").Value = 8545.5
$ws.Range("M135").Value = -6010.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 13234.75
$ws.Range("I43").Value = 10499.5
$ws.Range("K43").Value = 10499.5
$ws.Range("M43").Value = -10348.5

$ws.Range("H102").Value = 3115.75
$ws.Range("I102").Value = 3222.682
$ws.Range("K102").Value = 3222.682
$ws.Range("M102").Value = -1600.682

$ws.Range("H132").Value = 88844.664
$ws.Range("I132").Value = 88844.664
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 266533.992
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -264003.992
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 16061.875
$ws.Range("I68").Value = 4082.5
$ws.Range("K68").Value = 4082.5
$ws.Range("M68").Value = -3333.5

$ws.Range("H71").Value = 16061.875
$ws.Range("I71").Value = 4082.5
$ws.Range("K71").Value = 20412.5
$ws.Range("M71").Value = -16668.5

$ws.Range("H93").Value = 31251468
$ws.Range("I93").Value = 1381.25
$ws.Range("K93").Value = 1381.25
$ws.Range("M93").Value = -133.25

$ws.Range("H100").Value = 2215.6875
$ws.Range("I100").Value = 2096.7334
$ws.Range("K100").Value = 2096.7334
$ws.Range("M100").Value = -1555.7334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = 0

$ws.Range("H62").Value = 8474.556
$ws.Range("I62").Value = 1901
$ws.Range("K62").Value = 1901
$ws.Range("M62").Value = -1277

$ws.Range("H65").Value = 8474.556
$ws.Range("I65").Value = 1901
$ws.Range("K65").Value = 9505
$ws.Range("M65").Value = -6385

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = 0

$ws.Range("H136").Value = 528685.8
$ws.Range("I136").Value = 557946.25
$ws.Range("J136").Value = 1998
$ws.Range("K136").Value = 1673838.75
$ws.Range("L136").Value = 5994
$ws.Range("M136").Value = -1671288.75
$ws.Range("N136").Value = -11094
